$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 185
$ws1.Range("F3").Value = 965
$ws1.Range("F4").Value = 1119
$ws1.Range("F5").Value = 1566
$ws1.Range("F6").Value = 348
$ws1.Range("F7").Value = 721
$ws1.Range("F8").Value = 13018
$ws1.Range("F9").Value = 2255
$ws1.Range("F11").Value = 293
$ws1.Range("F12").Value = 53389
$ws1.Range("F13").Value = 1285
$ws1.Range("F14").Value = 292
$ws1.Range("F15").Value = 295
$ws1.Range("F16").Value = 847
$ws1.Range("F17").Value = 697
$ws1.Range("F20").Value = 835
$ws1.Range("F21").Value = 5017
$ws1.Range("F22").Value = 1226
$ws1.Range("F23").Value = 919
$ws1.Range("F28").Value = 1172
$ws1.Range("F30").Value = 20
$ws1.Range("F32").Value = 317
$ws1.Range("F35").Value = 55
$ws1.Range("F37").Value = 4647
$ws1.Range("F39").Value = 4701
$ws1.Range("F40").Value = 5653
$ws1.Range("F45").Value = 396
$ws1.Range("F46").Value = 93
$ws1.Range("F47").Value = 63
$ws1.Range("F49").Value = 166

# Sheet "演出" (sheet2) — column F update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 1095

# Sheet "本地生活" (sheet3) — column F updates
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 526
$ws3.Range("F4").Value = 132
$ws3.Range("F5").Value = 27

# Sheet "全部类型" (sheet4) — column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 526
$ws4.Range("F4").Value = 185
$ws4.Range("F5").Value = 965
$ws4.Range("F6").Value = 1119
$ws4.Range("F7").Value = 721
$ws4.Range("F8").Value = 13018
$ws4.Range("F9").Value = 13018
$ws4.Range("F10").Value = 2255
$ws4.Range("F12").Value = 292
$ws4.Range("F13").Value = 847
$ws4.Range("F14").Value = 697
$ws4.Range("F17").Value = 835
$ws4.Range("F19").Value = 5017
$ws4.Range("F20").Value = 1227
$ws4.Range("F21").Value = 27
$ws4.Range("F26").Value = 1172
$ws4.Range("F29").Value = 20
$ws4.Range("F32").Value = 317
$ws4.Range("F35").Value = 4647
$ws4.Range("F36").Value = 4701
$ws4.Range("F37").Value = 5653
$ws4.Range("F43").Value = 93
